$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 50
$ws.Range("B25").Value = "Update index.py"
$ws.Range("C25").Value = "riya-morankar"
$ws.Range("D25").Value = "N/A"
$ws.Range("E25").Value = "edit1 to main"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "2025-06-19"
